$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Tipo1 (F2) corrected from 21 to 22 ---
$ws.Range("F2").Value = 22

# --- New row 3: a second invoice line for the same customer/date,
#     mirroring row 2's layout but with invoice number A19262 and
#     Tipo1 (F) back at 21 ---
$ws.Range("A3").Value = "A19262"

# Copy the date cell from B2 so it lands as the same shared string
# (a literal "10/04/2024") instead of being re-parsed into a date serial.
$ws.Range("B2").Copy($ws.Range("B3"))

$ws.Range("C3").Value = "X3971208M"
$ws.Range("D3").Value = "PAPP ZSOLT"
$ws.Range("E3").Value = 65.52
$ws.Range("F3").Value = 21
$ws.Range("G3").Value = 13.76
$ws.Range("J3").Value = 12.25
$ws.Range("K3").Value = 10
$ws.Range("L3").Value = 1.23
$ws.Range("T3").Value = 20.83
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0
$ws.Range("AD3").Value = 344
$ws.Range("AF3").Value = 430000344

# --- Selection moves to F4 ---
[void]$ws.Range("F4").Select()

# --- Page setup: portrait, paper size 9 (A4) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
